$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.8000006675720215
$ws.Range("E2").Value = 214.7213788980116
$ws.Range("F2").Value = 0.005346381028066423
$ws.Range("G2").Value = 0.00492640138066684
$ws.Range("H2").Value = 0.00492640138066684
$ws.Range("I2").Value = 0.004602094923975882
$ws.Range("J2").Value = 0.004602094923975882
$ws.Range("K2").Value = 0.004602094923975882
$ws.Range("L2").Value = 0.004576315054447534
$ws.Range("M2").Value = 0.004575112486114983
$ws.Range("N2").Value = 0.004525511900407371
$ws.Range("O2").Value = 0.00440430735705288
$ws.Range("P2").Value = 0.00440430735705288
$ws.Range("Q2").Value = 0.00440430735705288
$ws.Range("R2").Value = 0.00432890496548281
$ws.Range("S2").Value = 0.004311766279417149
$ws.Range("T2").Value = 0.004295986399468312
$ws.Range("U2").Value = 0.004217292550986181
$ws.Range("V2").Value = 0.004217292550986181
$ws.Range("W2").Value = 0.00420888711460172
$ws.Range("X2").Value = 0.00420888711460172
$ws.Range("Y2").Value = 0.004185601927836483

$ws.Range("C3").Value = 0.8019962310791016
$ws.Range("E3").Value = 214.9350200352128
$ws.Range("F3").Value = 0.005406260595095675
$ws.Range("G3").Value = 0.004723112793375686
$ws.Range("H3").Value = 0.004723112793375686
$ws.Range("I3").Value = 0.004537026943229857
$ws.Range("J3").Value = 0.004537026943229857
$ws.Range("K3").Value = 0.004537026943229857
$ws.Range("L3").Value = 0.004537026943229857
$ws.Range("M3").Value = 0.004534997683126152
$ws.Range("N3").Value = 0.004503787047543423
$ws.Range("O3").Value = 0.004456764632231795
$ws.Range("P3").Value = 0.004438596941224078
$ws.Range("Q3").Value = 0.00440456944059498
$ws.Range("R3").Value = 0.004390166852866876
$ws.Range("S3").Value = 0.004309459224301395
$ws.Range("T3").Value = 0.004305706287762442
$ws.Range("U3").Value = 0.004286081872980953
$ws.Range("V3").Value = 0.004286081872980953
$ws.Range("W3").Value = 0.004241507753862991
$ws.Range("X3").Value = 0.004222944759700292
$ws.Range("Y3").Value = 0.004189766472421301

$ws.Range("C4").Value = 0.9010047912597656
$ws.Range("E4").Value = 214.5405462753642
$ws.Range("F4").Value = 0.005406260595095675
$ws.Range("G4").Value = 0.004934106956476759
$ws.Range("H4").Value = 0.004760832846073597
$ws.Range("I4").Value = 0.004667912507990768
$ws.Range("J4").Value = 0.004533212240490882
$ws.Range("K4").Value = 0.004455080191979171
$ws.Range("L4").Value = 0.004401794125058489
$ws.Range("M4").Value = 0.004380668526284448
$ws.Range("N4").Value = 0.004380668526284448
$ws.Range("O4").Value = 0.004380668526284448
$ws.Range("P4").Value = 0.004380668526284448
$ws.Range("Q4").Value = 0.004380668526284448
$ws.Range("R4").Value = 0.004290161955327986
$ws.Range("S4").Value = 0.00428819424405716
$ws.Range("T4").Value = 0.00428819424405716
$ws.Range("U4").Value = 0.004265671670227018
$ws.Range("V4").Value = 0.004214118961672974
$ws.Range("W4").Value = 0.004206368644095846
$ws.Range("X4").Value = 0.004185103729166167
$ws.Range("Y4").Value = 0.004182076925445696

$ws.Range("C5").Value = 0.7700023651123047
$ws.Range("E5").Value = 208.742485547662
$ws.Range("F5").Value = 0.005406260595095675
$ws.Range("G5").Value = 0.004952794870520366
$ws.Range("H5").Value = 0.004737167976803455
$ws.Range("I5").Value = 0.004650015505544784
$ws.Range("J5").Value = 0.004465853983348251
$ws.Range("K5").Value = 0.004463074221171216
$ws.Range("L5").Value = 0.004463074221171216
$ws.Range("M5").Value = 0.004277794001772648
$ws.Range("N5").Value = 0.004277794001772648
$ws.Range("O5").Value = 0.004277794001772648
$ws.Range("P5").Value = 0.004264414629703989
$ws.Range("Q5").Value = 0.004216604915438427
$ws.Range("R5").Value = 0.004216604915438427
$ws.Range("S5").Value = 0.004193022067826964
$ws.Range("T5").Value = 0.004163245994325718
$ws.Range("U5").Value = 0.00414873103165923
$ws.Range("V5").Value = 0.004119497605891216
$ws.Range("W5").Value = 0.004107912204862087
$ws.Range("X5").Value = 0.004101156945179413
$ws.Range("Y5").Value = 0.004069054299174696

$ws.Range("C6").Value = 0.8170065879821777
$ws.Range("E6").Value = 208.9947254736744
$ws.Range("F6").Value = 0.00532778429072544
$ws.Range("G6").Value = 0.004964024018688902
$ws.Range("H6").Value = 0.004704487760723903
$ws.Range("I6").Value = 0.004704487760723903
$ws.Range("J6").Value = 0.004659651212391287
$ws.Range("K6").Value = 0.004532514880877527
$ws.Range("L6").Value = 0.004466151366116666
$ws.Range("M6").Value = 0.004416028305618454
$ws.Range("N6").Value = 0.004315255861456995
$ws.Range("O6").Value = 0.00430188793347865
$ws.Range("P6").Value = 0.004206338340073381
$ws.Range("Q6").Value = 0.004206338340073381
$ws.Range("R6").Value = 0.004206338340073381
$ws.Range("S6").Value = 0.004161320120916609
$ws.Range("T6").Value = 0.004149372463692089
$ws.Range("U6").Value = 0.004136789144080959
$ws.Range("V6").Value = 0.004130591296758622
$ws.Range("W6").Value = 0.004118753219521267
$ws.Range("X6").Value = 0.004080085205518378
$ws.Range("Y6").Value = 0.00407397125679677

$ws.Range("C7").Value = 0.8060007095336914
$ws.Range("E7").Value = 206.7685864213654
$ws.Range("F7").Value = 0.005351429160723002
$ws.Range("G7").Value = 0.004897368424727533
$ws.Range("H7").Value = 0.004786342656350622
$ws.Range("I7").Value = 0.004644281545994164
$ws.Range("J7").Value = 0.004457881305069886
$ws.Range("K7").Value = 0.004457881305069886
$ws.Range("L7").Value = 0.004417133827596505
$ws.Range("M7").Value = 0.00430536331250398
$ws.Range("N7").Value = 0.00430536331250398
$ws.Range("O7").Value = 0.00430536331250398
$ws.Range("P7").Value = 0.004302434824660886
$ws.Range("Q7").Value = 0.004242886553701251
$ws.Range("R7").Value = 0.004145423972174123
$ws.Range("S7").Value = 0.004145423972174123
$ws.Range("T7").Value = 0.004137054187417091
$ws.Range("U7").Value = 0.004122899593527235
$ws.Range("V7").Value = 0.004102884122555275
$ws.Range("W7").Value = 0.004075172449718533
$ws.Range("X7").Value = 0.004043355312240814
$ws.Range("Y7").Value = 0.004030576733359948

$ws.Range("C8").Value = 0.826998233795166
$ws.Range("E8").Value = 209.8327443252747
$ws.Range("F8").Value = 0.005375769854281581
$ws.Range("G8").Value = 0.005039727683530361
$ws.Range("H8").Value = 0.004762905667361456
$ws.Range("I8").Value = 0.004635015773427452
$ws.Range("J8").Value = 0.004586353351731633
$ws.Range("K8").Value = 0.004586353351731633
$ws.Range("L8").Value = 0.004423979765568672
$ws.Range("M8").Value = 0.004297330577873542
$ws.Range("N8").Value = 0.004297330577873542
$ws.Range("O8").Value = 0.004297330577873542
$ws.Range("P8").Value = 0.004297330577873542
$ws.Range("Q8").Value = 0.004192647190263873
$ws.Range("R8").Value = 0.004192647190263873
$ws.Range("S8").Value = 0.004192647190263873
$ws.Range("T8").Value = 0.004192647190263873
$ws.Range("U8").Value = 0.004192647190263873
$ws.Range("V8").Value = 0.004129434550926653
$ws.Range("W8").Value = 0.004129434550926653
$ws.Range("X8").Value = 0.004108471072994321
$ws.Range("Y8").Value = 0.004090306906925432

$ws.Range("C9").Value = 0.7749979496002197
$ws.Range("E9").Value = 207.5377270594709
$ws.Range("F9").Value = 0.005232389430979058
$ws.Range("G9").Value = 0.004898309901870878
$ws.Range("H9").Value = 0.004854656145856027
$ws.Range("I9").Value = 0.00477115642799745
$ws.Range("J9").Value = 0.004511315616534178
$ws.Range("K9").Value = 0.004511315616534178
$ws.Range("L9").Value = 0.004511315616534178
$ws.Range("M9").Value = 0.004400354870813493
$ws.Range("N9").Value = 0.004357595374129738
$ws.Range("O9").Value = 0.004357595374129738
$ws.Range("P9").Value = 0.004341852558388379
$ws.Range("Q9").Value = 0.004217561468763999
$ws.Range("R9").Value = 0.004217561468763999
$ws.Range("S9").Value = 0.004162682399043533
$ws.Range("T9").Value = 0.004102329629303529
$ws.Range("U9").Value = 0.00408203763875929
$ws.Range("V9").Value = 0.004078823614020622
$ws.Range("W9").Value = 0.00406267642217193
$ws.Range("X9").Value = 0.004053013571449988
$ws.Range("Y9").Value = 0.004045569728254792

$ws.Range("C10").Value = 0.7310268878936768
$ws.Range("E10").Value = 214.7453286753334
$ws.Range("F10").Value = 0.005242365997607147
$ws.Range("G10").Value = 0.004917740760704342
$ws.Range("H10").Value = 0.004797743374100581
$ws.Range("I10").Value = 0.00464129355178569
$ws.Range("J10").Value = 0.00464129355178569
$ws.Range("K10").Value = 0.004536494475175598
$ws.Range("L10").Value = 0.004536494475175598
$ws.Range("M10").Value = 0.004485678602205299
$ws.Range("N10").Value = 0.004485678602205299
$ws.Range("O10").Value = 0.004395483629402811
$ws.Range("P10").Value = 0.004395483629402811
$ws.Range("Q10").Value = 0.004350592257498015
$ws.Range("R10").Value = 0.004321492851981777
$ws.Range("S10").Value = 0.004251742144040701
$ws.Range("T10").Value = 0.004251742144040701
$ws.Range("U10").Value = 0.004251742144040701
$ws.Range("V10").Value = 0.004242247807156977
$ws.Range("W10").Value = 0.004241557047214308
$ws.Range("X10").Value = 0.004224976585716723
$ws.Range("Y10").Value = 0.004186068785094218

$ws.Range("C11").Value = 0.7109959125518799
$ws.Range("E11").Value = 211.147855791025
$ws.Range("F11").Value = 0.005381263431984753
$ws.Range("G11").Value = 0.004922725300291017
$ws.Range("H11").Value = 0.004666652334736427
$ws.Range("I11").Value = 0.004666652334736427
$ws.Range("J11").Value = 0.004631370198230247
$ws.Range("K11").Value = 0.004610700792548699
$ws.Range("L11").Value = 0.004563682663676948
$ws.Range("M11").Value = 0.004563682663676948
$ws.Range("N11").Value = 0.004340619437376556
$ws.Range("O11").Value = 0.004340619437376556
$ws.Range("P11").Value = 0.004340619437376556
$ws.Range("Q11").Value = 0.0042020077078082
$ws.Range("R11").Value = 0.0042020077078082
$ws.Range("S11").Value = 0.0042020077078082
$ws.Range("T11").Value = 0.004198102479953639
$ws.Range("U11").Value = 0.004141757046547237
$ws.Range("V11").Value = 0.004141533744451154
$ws.Range("W11").Value = 0.004118884226401016
$ws.Range("X11").Value = 0.004118884226401016
$ws.Range("Y11").Value = 0.004115942608012182
